$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: "CON2-7" silkscreen reference became "CON2,3" now that CON4-7 ---
# --- gets its own line item (row 24) below.                               ---
$ws.Range("A9").Value = "CON2,3"

# --- New row 24: JST 2-pin connector (the "fixed qty JST 2 pin conn") ---
$ws.Range("A24").Value = "CON4-7"
$ws.Range("B24").Value = "CONN HEADER GH SIDE 2POS 1.25MM"
$ws.Range("C24").Value = "SM02B-GHS-TB(LF)(SN)"
$ws.Range("D24").Value = "http://www.jst-mfg.com/product/pdf/eng/eGH.pdf"
$ws.Range("E24").Value = 0.46
$ws.Range("F24").Value = 4
$ws.Range("G24").Formula = "=F24*E24"

# Match the plain bordered style used by the rest of the parts table.
$ws.Range("A24:G24").Borders.LineStyle = 1

# Wire up the live hyperlink, then reuse the hyperlink-cell formatting
# (border + hyperlink font) from the existing Datasheet column so the
# cell style matches the rest of the table instead of Excel's default
# "newly-added hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("D24"), "http://www.jst-mfg.com/product/pdf/eng/eGH.pdf") | Out-Null
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4122) | Out-Null

# Move the active selection the way the author's session ended up.
$ws.Range("J25").Select() | Out-Null
